$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original row 2 (sending cluster "MuSCs" -> "a" -> "Mc2r" -> "FAPs") needs to
# move down to row 3, and a new row 2 is inserted above it for sending cluster "ECs"
# (same ligand/receptor/target cluster triple), with the total-expression-derived
# specificity columns (I, J, S, T) recomputed across both rows so they sum to 1.

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Copy current row 2 values down into row 3 (plain copy, no style carried along).
foreach ($col in $cols) {
    $srcAddr = $col + "2"
    $dstAddr = $col + "3"
    $ws.Range($dstAddr).Value = $ws.Range($srcAddr).Value2
}

# New row 2: "ECs" sending cluster row with its own computed values.
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "a"
$ws.Range("C2").Value = "Mc2r"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.045339
$ws.Range("H2").Value = 0.136017
$ws.Range("I2").Value = 0.1740293637846656
$ws.Range("J2").Value = 0.1740293637846656
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04473366666666666
$ws.Range("N2").Value = 0.134201
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.002028179713
$ws.Range("R2").Value = 0.018253617417
$ws.Range("S2").Value = 0.1740293637846656
$ws.Range("T2").Value = 0.1740293637846656

# Row 3 ("MuSCs") keeps its original literal values except the derived-specificity
# columns, which are recomputed now that "ECs" shares the total.
$ws.Range("I3").Value = 0.8259706362153344
$ws.Range("J3").Value = 0.8259706362153345
$ws.Range("S3").Value = 0.8259706362153344
$ws.Range("T3").Value = 0.8259706362153345
